$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.420.32'
$ws.Range("E2").Value = '  -4.05%  '

# Row 3
$ws.Range("D3").Value = '1.770.73'
$ws.Range("E3").Value = '  -3.24%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("E5").Value = '  +0.20%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.20%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4288'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.61%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3650'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07181'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.51%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8490'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.62%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.37'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.34%  '

# Row 12
$ws.Range("D12").Value = '1.796.26'
$ws.Range("E12").Value = '  -4.72%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.441'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.79%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.243'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.84%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06896'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.31%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.17%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.43%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008684'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.89%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.14%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.33%  '

# Row 21
$ws.Range("D21").Value = '26.424.79'
$ws.Range("E21").Value = '  -4.83%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.123'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.59%  '

# Row 24
$ws.Range("D24").Value = '1.991.23'
$ws.Range("E24").Value = '  -5.64%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.05%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.860'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.22%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.48%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.083'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.51%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.01%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.752'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.45%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08947'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.18%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7251'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.17%  '

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.116'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.72%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.336'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.81%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.001'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.742'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.93%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.079'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.68%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05166'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.82%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01887'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.04%  '

# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1614'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.48%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4919'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.19%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.575'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.73%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.284'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.31%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.037'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.74'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.92%  '

# Row 46
$ws.Range("E46").Value = '  -2.81%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.20%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06193'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.22%  '

# Row 49
$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4472'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.78%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.594'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.52%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.745'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.53%  '
